$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 43876
$ws.Range("B47").NumberFormat = "d-mmm-yy"
$ws.Range("C47").Value = "Surah Al Nisa, 36 - 41"
$d47 = @'
h1: Winners never Quit and Quitters never win
img: 46.jpg
p.caption.text-center: This is me and my brother in June 2017, before I got admitted (Aug 2017) to the psychiatry ward
p: I was thrashed brutally by my close family members yesterday over writing these blogs everyday and organizing monthly tech meetups. I was asked if I am taking my psychiatry medications or not. Weather I am living a sane life or not. I was also told, “<b>It is wrong to bring our daughters and sisters to tech meetups where I have also invited my close friends</b>”. Some of these relatives are rich and have accumulated good wealth doing business abroad. I was also asked if the books I am reading are complete enough to dictate the true picture of Quran.
p: I will try to write my answer to these queries, because I am a Telecom Engineer, who can make websites and write blogs everyday. I will try to act proud and blunt because I really do not care what people think of a person who writes everyday about <span class=lavendar>“Quran”</span>. 
quote: (O Muhammad) Say, "I only advise you of one [thing] -  <span class=lavendar>that you stand for Allah, [seeking truth] </span>in pairs and individually, and then give thought." <span class=lavendar>There is not in your companion any madness</span>. He is only a warner to you before a severe punishment. <br> - Surah Saba verse 46
h3: How did I become a psychiatry patient?
p: I started taking pictures everyday in <b>2017</b>, <span class=note>like the one above</span>. It went on for 60 days which lead me believing I can give Friday Sermon in our local mosque. I showed the intent at my profession to express myself. 
p: Because it is unorthodox to give sermons for an engineer, I was made admitted in a psychiatry ward for 1 month. It is a prison cell where you can not move out. I was labelled a “<b>Mood Disorder</b>” patient.
p: Back then I did not have a protection wall around me, which made it relatively easier for people to get me to the prison.
h3: 2018, My second time in psychiatry ward.
p: I tried to cope with the high dose epilepsy medication. I managed to taper it off, to stand back on my feet and I again wanted to give a sermon and showed the intent. I was again considered a psychiatry patient. My mother and brother started giving me the medicine without my knowledge. As soon as I felt the dose of medication in my urine, I confirmed it from them. I found out they secretly were poring high dose epilepsy medicine into my food for 3 - 4 days. 
p: I surrendered and went on my own to the same Psychiatry ward for 2 months. My wife had not yet moved to Pakistan.
h3: Building of protection wall
h5: 1. Wife
p: I married to a Norwegian girl. She is a strong, bold, blunt and ruthless girl. She fought with the head of psychiatry doctor and lowered him to the ground. I was discharged from the hospital with positive remarks. I joined my profession back. I worked hard and got an <b>out-standing report</b> for performing an excellent job from my boss.
h5: 2. Success at University
p: I started doing Masters. I launched <a href=https://www.zakatlists.com/school target=_blank>zakatlists.com/school</a>. My university class fellows are using it. I have topped 2 subjects out of 3 in the first semester of university. I do not know yet the result of the 3rd subject.
h5: 3. Friends over Quran
p: I joined a circle of well educated friends who read Quran everyday in the morning and read the translation. I write this blog daily where I try to write about my feelings and experiences towards Quran. The ways and means through which it has impacted my life.
p: I have my wife and these friends made over Quran, beside me. If my blog writing is abnormal, I do not believe, labelling my wife and close friends ‘abnormal’ is possible. I can not be victimised any more because of my beliefs and actions.
p: The building of this wall has made it possible for me that I keep writing a blog daily. I also arrange tech Meetups every month. I will not stop. 
h3: Conclusion
p: It is not an easy job to arrange a tech meetup every month. I am paying Rs 20000 every month out of my own pocket, because I am not a quitter. I will not lower my voice. I will not stop writing. I will not stop reading Quran. Because I am a proud guy, who inspite of getting labelled a ‘psychiatry patient’ will not get discouraged.
p: These blogs under <a href=https://www.zakatlists.com target=_blank>zakatlists.com</a> might turn out making Pakistan better. It might help Imran Khan. It might help our youth love our mosques. <b>And above all, it might unite us; from Indonesia to Palestine</b>. Who knows, it might become something bigger. 
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>. Subscribe below to receive updates.
'@
$ws.Range("D47").Value = $d47
$ws.Range("E47").Value = "Qasim Ali"
$ws.Range("F47").Value = "Techshek, zakatlists, blogging, Quran daily, Will I stop?"

$ws.Rows.Item(47).RowHeight = 409.6

$ws.Range("F47").Select()
